# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71536.68398845789

# Row 3
$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 250.0739139718791
